$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 116
$ws.Range("H116").Value = 3976.5833
$ws.Range("I116").Value = 3116.3076
$ws.Range("J116").Value = 4993.273
$ws.Range("K116").Value = 3116.3076
$ws.Range("L116").Value = 4993.273
$ws.Range("M116").Value = 325.6923999999999
$ws.Range("N116").Value = -11877.273

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 3280
$ws.Range("I63").Value = 2960
$ws.Range("K63").Value = 2960
$ws.Range("M63").Value = -2274

# Row 66
$ws.Range("H66").Value = 3280
$ws.Range("I66").Value = 2960
$ws.Range("K66").Value = 14800
$ws.Range("M66").Value = -11368

# Row 74
$ws.Range("H74").Value = 865.6
$ws.Range("I74").Value = 824.2353000000001
$ws.Range("J74").Value = 1100
$ws.Range("K74").Value = 824.2353000000001
$ws.Range("L74").Value = 1100
$ws.Range("M74").Value = 49.76469999999995
$ws.Range("N74").Value = -2848

# Row 77
$ws.Range("H77").Value = 865.6
$ws.Range("I77").Value = 824.2353000000001
$ws.Range("J77").Value = 1100
$ws.Range("K77").Value = 4121.1765
$ws.Range("L77").Value = 5500
$ws.Range("M77").Value = 246.8234999999995
$ws.Range("N77").Value = -14236

# Row 102
$ws.Range("H102").Value = 5833.1665
$ws.Range("I102").Value = 4999.8
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 4999.8
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -3377.8
$ws.Range("N102").Value = -13244

# Row 122
$ws.Range("H122").Value = 2402.4
$ws.Range("I122").Value = 1802.8334
$ws.Range("J122").Value = 3301.75
$ws.Range("K122").Value = 5408.5002
$ws.Range("L122").Value = 9905.25
$ws.Range("M122").Value = -2958.5002
$ws.Range("N122").Value = -14805.25

# Row 132
$ws.Range("H132").Value = 21279164
$ws.Range("I132").Value = 30304680
$ws.Range("J132").Value = 4732.357
$ws.Range("K132").Value = 90914040
$ws.Range("L132").Value = 14197.071
$ws.Range("M132").Value = -90911510
$ws.Range("N132").Value = -19257.071

$ws = $wb.Worksheets.Item("BSM")
# Row 45
$ws.Range("H45").Value = 26731.666
$ws.Range("J45").Value = 26731.666
$ws.Range("L45").Value = 26731.666
$ws.Range("N45").Value = -28347.666

# Row 134
$ws.Range("H134").Value = 1847.3243
$ws.Range("I134").Value = 1098.1613
$ws.Range("J134").Value = 5718
$ws.Range("K134").Value = 3294.4839
$ws.Range("L134").Value = 17154
$ws.Range("M134").Value = -759.4839000000002
$ws.Range("N134").Value = -22224

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2445.392
$ws.Range("I31").Value = 1902.6136
$ws.Range("J31").Value = 5857.143
$ws.Range("K31").Value = 1902.6136
$ws.Range("L31").Value = 5857.143
$ws.Range("M31").Value = -1607.6136
$ws.Range("N31").Value = -6447.143

# Row 34
$ws.Range("H34").Value = 2445.392
$ws.Range("I34").Value = 1902.6136
$ws.Range("J34").Value = 5857.143
$ws.Range("K34").Value = 1902.6136
$ws.Range("L34").Value = 5857.143
$ws.Range("M34").Value = -1700.6136
$ws.Range("N34").Value = -6261.143

# Row 47
$ws.Range("H47").Value = 30000
$ws.Range("J47").Value = 30000
$ws.Range("L47").Value = 30000
$ws.Range("N47").Value = -31132

# Row 62
$ws.Range("H62").Value = 3489.9167
$ws.Range("I62").Value = 2295.8
$ws.Range("J62").Value = 4342.857
$ws.Range("K62").Value = 2295.8
$ws.Range("L62").Value = 4342.857
$ws.Range("M62").Value = -1671.8
$ws.Range("N62").Value = -5590.857

# Row 65
$ws.Range("H65").Value = 3489.9167
$ws.Range("I65").Value = 2295.8
$ws.Range("J65").Value = 4342.857
$ws.Range("K65").Value = 11479
$ws.Range("L65").Value = 21714.285
$ws.Range("M65").Value = -8359
$ws.Range("N65").Value = -27954.285

# Row 134
$ws.Range("H134").Value = 2021.8422
$ws.Range("I134").Value = 978.1539
$ws.Range("J134").Value = 4283.1665
$ws.Range("K134").Value = 2934.4617
$ws.Range("L134").Value = 12849.4995
$ws.Range("M134").Value = -399.4616999999998
$ws.Range("N134").Value = -17919.4995

$ws = $wb.Worksheets.Item("CUL")
# Row 82
$ws.Range("H82").Value = 3060
$ws.Range("I82").Value = 300
$ws.Range("J82").Value = 3750
$ws.Range("K82").Value = 900
$ws.Range("L82").Value = 11250
$ws.Range("M82").Value = -494
$ws.Range("N82").Value = -12062

# Row 85
$ws.Range("H85").Value = 3060
$ws.Range("I85").Value = 300
$ws.Range("J85").Value = 3750
$ws.Range("K85").Value = 900
$ws.Range("L85").Value = 11250
$ws.Range("M85").Value = 504
$ws.Range("N85").Value = -14058

# Row 131
$ws.Range("H131").Value = 1205.4237
$ws.Range("I131").Value = 2865.6
$ws.Range("J131").Value = 1051.7037
$ws.Range("K131").Value = 8596.799999999999
$ws.Range("L131").Value = 3155.1111
$ws.Range("M131").Value = -3556.799999999999
$ws.Range("N131").Value = -13235.1111

$ws = $wb.Worksheets.Item("GSM")
# Row 137
$ws.Range("H137").Value = 29623.75
$ws.Range("J137").Value = 29623.75
$ws.Range("L137").Value = 29623.75
$ws.Range("N137").Value = -39823.75

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3400
$ws.Range("I7").Value = 1500
$ws.Range("J7").Value = 4666.6665
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 4666.6665
$ws.Range("M7").Value = -1388
$ws.Range("N7").Value = -4890.6665

# Row 55
$ws.Range("H55").Value = 777.5238000000001
$ws.Range("I55").Value = 218.54546
$ws.Range("J55").Value = 1392.4
$ws.Range("K55").Value = 218.54546
$ws.Range("L55").Value = 1392.4
$ws.Range("M55").Value = -45.54545999999999
$ws.Range("N55").Value = -1738.4

# Row 60
$ws.Range("H60").Value = 18040.666
$ws.Range("J60").Value = 18040.666
$ws.Range("L60").Value = 18040.666
$ws.Range("N60").Value = -19058.666

# Row 68
$ws.Range("H68").Value = 1433.2916
$ws.Range("I68").Value = 1022.6818
$ws.Range("J68").Value = 5950
$ws.Range("K68").Value = 1022.6818
$ws.Range("L68").Value = 5950
$ws.Range("M68").Value = -273.6818
$ws.Range("N68").Value = -7448

# Row 71
$ws.Range("H71").Value = 1433.2916
$ws.Range("I71").Value = 1022.6818
$ws.Range("J71").Value = 5950
$ws.Range("K71").Value = 5113.409
$ws.Range("L71").Value = 29750
$ws.Range("M71").Value = -1369.409
$ws.Range("N71").Value = -37238

# Row 122
$ws.Range("H122").Value = 3194.524
$ws.Range("I122").Value = 2248.8333
$ws.Range("J122").Value = 4455.4443
$ws.Range("K122").Value = 6746.499899999999
$ws.Range("L122").Value = 13366.3329
$ws.Range("M122").Value = -4296.499899999999
$ws.Range("N122").Value = -18266.3329

# Row 126
$ws.Range("H126").Value = 3400
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 4666.6665
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 13999.9995
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -18939.9995

# Row 136
$ws.Range("H136").Value = 1945.9474
$ws.Range("I136").Value = 1469.258
$ws.Range("J136").Value = 4057
$ws.Range("K136").Value = 4407.774
$ws.Range("L136").Value = 12171
$ws.Range("M136").Value = -1857.774
$ws.Range("N136").Value = -17271

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 3450382.2
$ws.Range("I126").Value = 1799
$ws.Range("K126").Value = 5397
$ws.Range("M126").Value = -2927

# Row 132
$ws.Range("H132").Value = 3075.831
$ws.Range("I132").Value = 1226.8928
$ws.Range("J132").Value = 9978.532999999999
$ws.Range("K132").Value = 3680.6784
$ws.Range("L132").Value = 29935.599
$ws.Range("M132").Value = -1150.6784
$ws.Range("N132").Value = -34995.599

# Row 133
$ws.Range("H133").Value = 39630
$ws.Range("J133").Value = 39630
$ws.Range("L133").Value = 39630
$ws.Range("N133").Value = -49750

# Row 136
$ws.Range("H136").Value = 801.0606
$ws.Range("I136").Value = 541.4912
$ws.Range("J136").Value = 2445
$ws.Range("K136").Value = 1624.4736
$ws.Range("L136").Value = 7335
$ws.Range("M136").Value = 925.5263999999997
$ws.Range("N136").Value = -12435
